$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "68.256.60"
$ws.Range("E2").Value = "  -1.82%  "

# Row 3
Set-TextValue "D3" "3.916.96"
$ws.Range("E3").Value = "  -1.65%  "

# Row 4
$ws.Range("E4").Value = "  +0.34%  "

# Row 5
Set-TextValue "D5" "484.99"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
Set-TextValue "D6" "146.28"
$ws.Range("E6").Value = "  -2.31%  "

# Row 7
Set-TextValue "D7" "0.622"
$ws.Range("E7").Value = "  -1.07%  "

# Row 8
Set-TextValue "D8" "0.997"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
Set-TextValue "D9" "0.737"
$ws.Range("E9").Value = "  -0.16%  "

# Row 10
$ws.Range("E10").Value = "  -1.16%  "

# Row 11
Set-TextValue "D11" "0.0000347"
$ws.Range("E11").Value = "  -3.55%  "

# Row 12
Set-TextValue "D12" "43.21"
$ws.Range("E12").Value = "  -0.40%  "

# Row 13
$ws.Range("E13").Value = "  +2.01%  "

# Row 14
Set-TextValue "D14" "4.539.05"
$ws.Range("E14").Value = "  -0.97%  "

# Row 15
Set-TextValue "D15" "3.917.25"
$ws.Range("E15").Value = "  -0.74%  "

# Row 16
Set-TextValue "D16" "14.24"
$ws.Range("E16").Value = "  -4.15%  "

# Row 17
$ws.Range("E17").Value = "  -0.79%  "

# Row 18
Set-TextValue "D18" "20.26"
$ws.Range("E18").Value = "  +0.90%  "

# Row 19
$ws.Range("E19").Value = "  -0.12%  "

# Row 20
Set-TextValue "D20" "68.280.66"
$ws.Range("E20").Value = "  -1.49%  "

# Row 21
Set-TextValue "D21" "431.11"
$ws.Range("E21").Value = "  -3.00%  "

# Row 22
Set-TextValue "D22" "3.53"
$ws.Range("E22").Value = "  +5.78%  "

# Row 23
Set-TextValue "D23" "15.09"
$ws.Range("E23").Value = "  +3.97%  "

# Row 24
Set-TextValue "D24" "89.22"
$ws.Range("E24").Value = "  +0.41%  "

# Row 25
Set-TextValue "D25" "11.76"
$ws.Range("E25").Value = "  +20.34%  "

# Row 26
Set-TextValue "D26" "3.72"
$ws.Range("E26").Value = "  +0.48%  "

# Row 27
$ws.Range("E27").Value = "  +10.11%  "

# Row 28
Set-TextValue "D28" "37.79"
$ws.Range("E28").Value = "  -3.24%  "

# Row 29
Set-TextValue "D29" "5.67"
$ws.Range("E29").Value = "  -1.79%  "

# Row 30
Set-TextValue "D30" "719.45"
$ws.Range("E30").Value = "  -1.86%  "

# Row 31
Set-TextValue "D31" "13.76"
$ws.Range("E31").Value = "  +3.46%  "

# Row 32
$ws.Range("E32").Value = "  +2.03%  "

# Row 33
Set-TextValue "D33" "2.94"
$ws.Range("E33").Value = "  +3.61%  "

# Row 34
Set-TextValue "D34" "6.20"
$ws.Range("E34").Value = "  +14.97%  "

# Row 35
Set-TextValue "D35" "0.0₃0889"
$ws.Range("E35").Value = "  +3.17%  "

# Row 36
Set-TextValue "D36" "41.70"
$ws.Range("E36").Value = "  -2.12%  "

# Row 37
Set-TextValue "D37" "61.13"
$ws.Range("E37").Value = "  +1.34%  "

# Row 38
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D38" "3.08"
$ws.Range("E38").Value = "  +19.37%  "

# Row 39
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D39" "0.403"
$ws.Range("E39").Value = "  +20.09%  "

# Row 40
Set-TextValue "D40" "0.997"
$ws.Range("E40").Value = "  -0.14%  "

# Row 41
$ws.Range("E41").Value = "  -4.50%  "

# Row 42
$ws.Range("E42").Value = "  +3.29%  "

# Row 43
$ws.Range("E43").Value = "  +1.65%  "

# Row 44
$ws.Range("E44").Value = "  +1.50%  "

# Row 45
$ws.Range("E45").Value = "  -0.19%  "

# Row 46
Set-TextValue "D46" "3.37"
$ws.Range("E46").Value = "  +3.41%  "

# Row 47
$ws.Range("E47").Value = "  +0.43%  "

# Row 48
Set-TextValue "D48" "3.44"
$ws.Range("E48").Value = "  +0.90%  "

# Row 49
$ws.Range("E49").Value = "  -2.38%  "

# Row 50
Set-TextValue "D50" "145.44"
$ws.Range("E50").Value = "  -2.16%  "

# Row 51
Set-TextValue "D51" "0.0₆0334"
$ws.Range("E51").Value = "  +27.90%  "
